$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value
$changes = @(
    @("D2", "43.429.58"),
    @("E2", "  -4.16%  "),
    @("D3", "2.237.74"),
    @("E3", "  -5.30%  "),
    @("D4", "1.01"),
    @("E4", "  +0.27%  "),
    @("D5", "319.61"),
    @("E5", "  +3.44%  "),
    @("D6", "100.59"),
    @("E6", "  -8.14%  "),
    @("E7", "  -7.20%  "),
    @("D9", "0.568"),
    @("E9", "  -7.84%  "),
    @("D10", "36.94"),
    @("E10", "  -10.38%  "),
    @("D11", "54.34"),
    @("E11", "  -1.83%  "),
    @("D12", "0.0824"),
    @("E12", "  -10.19%  "),
    @("D13", "7.74"),
    @("E13", "  -8.60%  "),
    @("E14", "  -2.93%  "),
    @("D15", "0.870"),
    @("E15", "  -11.58%  "),
    @("D16", "2.577.53"),
    @("E16", "  -5.22%  "),
    @("D17", "14.22"),
    @("E17", "  -7.42%  "),
    @("D18", "2.243.41"),
    @("E18", "  -4.96%  "),
    @("D19", "43.410.59"),
    @("E19", "  -4.14%  "),
    @("D20", "14.41"),
    @("E20", "  +6.70%  "),
    @("D21", "0.0₃0974"),
    @("E21", "  -8.66%  "),
    @("D22", "6.55"),
    @("E22", "  -10.34%  "),
    @("D23", "65.66"),
    @("E23", "  -10.41%  "),
    @("D24", "3.21"),
    @("E24", "  -6.76%  "),
    @("D25", "236.44"),
    @("E25", "  -8.74%  "),
    @("D26", "2.18"),
    @("E26", "  -5.66%  "),
    @("E27", "  -0.08%  "),
    @("D28", "10.25"),
    @("E28", "  -7.88%  "),
    @("E29", "  -7.14%  "),
    @("D30", "6.47"),
    @("E30", "  -12.03%  "),
    @("D31", "0.0893"),
    @("E31", "  -7.67%  "),
    @("D32", "20.69"),
    @("E32", "  -7.70%  "),
    @("D33", "34.26"),
    @("E33", "  -11.01%  "),
    @("D34", "158.63"),
    @("E34", "  -7.00%  "),
    @("D35", "2.77"),
    @("E35", "  -5.42%  "),
    @("D36", "3.35"),
    @("E36", "  +13.17%  "),
    @("B37", "ARBITRUM"),
    @("C37", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D37", "2.04"),
    @("E37", "  +17.33%  "),
    @("B38", "Stellar"),
    @("C38", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("D38", "0.123"),
    @("E38", "  -6.28%  "),
    @("D39", "4.53"),
    @("E39", "  -6.70%  "),
    @("D40", "0.106"),
    @("E40", "  -7.03%  "),
    @("D41", "3.65"),
    @("E41", "  -7.12%  "),
    @("D42", "0.0327"),
    @("E42", "  -8.50%  "),
    @("E43", "  +0.41%  "),
    @("D44", "1.817.61"),
    @("E44", "  +9.71%  "),
    @("D45", "12.15"),
    @("E45", "  -5.28%  "),
    @("D46", "88.81"),
    @("E46", "  -10.46%  "),
    @("B47", "THORChain"),
    @("C47", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"),
    @("D47", "5.56"),
    @("E47", "  +0.76%  "),
    @("B48", "Algorand"),
    @("C48", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"),
    @("D48", "0.209"),
    @("E48", "  -10.20%  "),
    @("B49", "ordi"),
    @("C49", "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"),
    @("D49", "79.09"),
    @("E49", "  -3.74%  "),
    @("B50", "EnergySwap"),
    @("C50", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D50", "16.45"),
    @("E50", "  +64.03%  "),
    @("D51", "61.34"),
    @("E51", "  -12.27%  "),
)

foreach ($change in $changes) {
    $addr = $change[0]
    $value = $change[1]
    $cell = $ws.Range($addr)
    # Force the cell to Text format first so Excel does not reinterpret
    # numeric-looking strings (e.g. "1.01") as actual numbers; this keeps
    # the stored cell type as a string, matching the source data which is
    # entirely textual (inline strings in the original workbook).
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
